$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Single-cell changes in column F (rows that keep their row number) ---
$ws.Range("F2").Value = 18.03
$ws.Range("F3").Value = ""
$ws.Range("F4").Value = ""
$ws.Range("F11").Value = 17.65
$ws.Range("F13").Value = ""
$ws.Range("F21").Value = 16.58
$ws.Range("F25").Value = ""

# --- Delete entire rows for "RM 232" (row 26) and "SC 92" (originally row 28,
#     becomes row 27 after the first deletion shifts rows up) ---
$ws.Rows.Item(26).Delete()
$ws.Rows.Item(27).Delete()

# --- After the deletions, remaining rows shifted up. Apply remaining value
#     changes to "SC 119" (now row 29) and "SC 232" (now row 33) ---
$ws.Range("D29").Value = ""
$ws.Range("D33").Value = -14.1
$ws.Range("F33").Value = 17.53
